$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 570, shifting existing rows 570:682 down to 571:683
$ws.Rows.Item(570).Insert()

# Populate the new row 570 with values (same group values as neighboring rows, new price data)
$ws.Cells.Item(570, 1).Value = 11
$ws.Cells.Item(570, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(570, 3).Value = "Bíobío"
$ws.Cells.Item(570, 4).Value = 45244
$ws.Cells.Item(570, 5).Value = 8
$ws.Cells.Item(570, 6).Value = 100112006
$ws.Cells.Item(570, 7).Value = "Repollo"
$ws.Cells.Item(570, 8).Value = "Crespo record"
$ws.Cells.Item(570, 9).Value = "Primera"
$ws.Cells.Item(570, 10).Value = 900
$ws.Cells.Item(570, 11).Value = 1000
$ws.Cells.Item(570, 12).Value = 1000
$ws.Cells.Item(570, 13).Value = 1000
$ws.Cells.Item(570, 14).Value = "`$/unidad"
$ws.Cells.Item(570, 15).Value = "Región Metropolitana"
$ws.Cells.Item(570, 16).Value = 1000
$ws.Cells.Item(570, 17).Value = 1
$ws.Cells.Item(570, 18).Value = "Hortaliza"
